$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.07164033333333
$ws.Range("H2").Value = 45.214921
$ws.Range("I2").Value = 0.04362521826456533
$ws.Range("J2").Value = 0.04362521826456534
$ws.Range("M2").Value = 12.23245266666667
$ws.Range("N2").Value = 36.697358
$ws.Range("O2").Value = 0.2605176191876535
$ws.Range("P2").Value = 0.2605176191876535
$ws.Range("Q2").Value = 184.3631269865242
$ws.Range("R2").Value = 1659.268142878718
$ws.Range("S2").Value = 0.0113651379988263
$ws.Range("T2").Value = 0.0113651379988263
$ws.Range("G3").Value = 15.07164033333333
$ws.Range("H3").Value = 45.214921
$ws.Range("I3").Value = 0.04362521826456533
$ws.Range("J3").Value = 0.04362521826456534
$ws.Range("O3").Value = 0.2463118768921188
$ws.Range("P3").Value = 0.2463118768921188
$ws.Range("Q3").Value = 174.3100062842234
$ws.Range("R3").Value = 1568.790056558011
$ws.Range("S3").Value = 0.01074540939057343
$ws.Range("T3").Value = 0.01074540939057343
$ws.Range("G4").Value = 15.07164033333333
$ws.Range("H4").Value = 45.214921
$ws.Range("I4").Value = 0.04362521826456533
$ws.Range("J4").Value = 0.04362521826456534
$ws.Range("M4").Value = 20.18080466666667
$ws.Range("N4").Value = 60.54241400000001
$ws.Range("O4").Value = 0.4297956696270414
$ws.Range("P4").Value = 0.4297956696270414
$ws.Range("Q4").Value = 304.157829573255
$ws.Range("R4").Value = 2737.420466159294
$ws.Range("S4").Value = 0.01874992989664469
$ws.Range("T4").Value = 0.01874992989664469
$ws.Range("G5").Value = 15.07164033333333
$ws.Range("H5").Value = 45.214921
$ws.Range("I5").Value = 0.04362521826456533
$ws.Range("J5").Value = 0.04362521826456534
$ws.Range("M5").Value = 2.975728333333333
$ws.Range("N5").Value = 8.927185
$ws.Range("O5").Value = 0.06337483429318623
$ws.Range("P5").Value = 0.06337483429318624
$ws.Range("Q5").Value = 44.84910716970944
$ws.Range("R5").Value = 403.641964527385
$ws.Range("S5").Value = 0.002764740978520909
$ws.Range("T5").Value = 0.00276474097852091
$ws.Range("I6").Value = 0.7276761574605194
$ws.Range("J6").Value = 0.7276761574605195
$ws.Range("M6").Value = 12.23245266666667
$ws.Range("N6").Value = 36.697358
$ws.Range("O6").Value = 0.2605176191876535
$ws.Range("P6").Value = 0.2605176191876535
$ws.Range("Q6").Value = 3075.208724673104
$ws.Range("R6").Value = 27676.87852205793
$ws.Range("S6").Value = 0.1895724600812346
$ws.Range("T6").Value = 0.1895724600812346
$ws.Range("I7").Value = 0.7276761574605194
$ws.Range("J7").Value = 0.7276761574605195
$ws.Range("O7").Value = 0.2463118768921188
$ws.Range("P7").Value = 0.2463118768921188
$ws.Range("S7").Value = 0.1792352801137455
$ws.Range("T7").Value = 0.1792352801137455
$ws.Range("I8").Value = 0.7276761574605194
$ws.Range("J8").Value = 0.7276761574605195
$ws.Range("M8").Value = 20.18080466666667
$ws.Range("N8").Value = 60.54241400000001
$ws.Range("O8").Value = 0.4297956696270414
$ws.Range("P8").Value = 0.4297956696270414
$ws.Range("Q8").Value = 5073.40500494807
$ws.Range("R8").Value = 45660.64504453263
$ws.Range("S8").Value = 0.3127520613673764
$ws.Range("T8").Value = 0.3127520613673764
$ws.Range("I9").Value = 0.7276761574605194
$ws.Range("J9").Value = 0.7276761574605195
$ws.Range("M9").Value = 2.975728333333333
$ws.Range("N9").Value = 8.927185
$ws.Range("O9").Value = 0.06337483429318623
$ws.Range("P9").Value = 0.06337483429318624
$ws.Range("Q9").Value = 748.0908352795005
$ws.Range("R9").Value = 6732.817517515505
$ws.Range("S9").Value = 0.04611635589816291
$ws.Range("T9").Value = 0.04611635589816292
$ws.Range("G10").Value = 24.41113566666667
$ws.Range("H10").Value = 73.233407
$ws.Range("I10").Value = 0.07065860768910215
$ws.Range("J10").Value = 0.07065860768910216
$ws.Range("M10").Value = 12.23245266666667
$ws.Range("N10").Value = 36.697358
$ws.Range("O10").Value = 0.2605176191876535
$ws.Range("P10").Value = 0.2605176191876535
$ws.Range("Q10").Value = 298.6080615820784
$ws.Range("R10").Value = 2687.472554238706
$ws.Range("S10").Value = 0.01840781225027932
$ws.Range("T10").Value = 0.01840781225027932
$ws.Range("G11").Value = 24.41113566666667
$ws.Range("H11").Value = 73.233407
$ws.Range("I11").Value = 0.07065860768910215
$ws.Range("J11").Value = 0.07065860768910216
$ws.Range("O11").Value = 0.2463118768921188
$ws.Range("P11").Value = 0.2463118768921188
$ws.Range("Q11").Value = 282.3252889103819
$ws.Range("R11").Value = 2540.927600193437
$ws.Range("S11").Value = 0.01740405427848665
$ws.Range("T11").Value = 0.01740405427848665
$ws.Range("G12").Value = 24.41113566666667
$ws.Range("H12").Value = 73.233407
$ws.Range("I12").Value = 0.07065860768910215
$ws.Range("J12").Value = 0.07065860768910216
$ws.Range("M12").Value = 20.18080466666667
$ws.Range("N12").Value = 60.54241400000001
$ws.Range("O12").Value = 0.4297956696270414
$ws.Range("P12").Value = 0.4297956696270414
$ws.Range("Q12").Value = 492.6363605804999
$ws.Range("R12").Value = 4433.727245224499
$ws.Range("S12").Value = 0.03036876360665208
$ws.Range("T12").Value = 0.03036876360665208
$ws.Range("G13").Value = 24.41113566666667
$ws.Range("H13").Value = 73.233407
$ws.Range("I13").Value = 0.07065860768910215
$ws.Range("J13").Value = 0.07065860768910216
$ws.Range("M13").Value = 2.975728333333333
$ws.Range("N13").Value = 8.927185
$ws.Range("O13").Value = 0.06337483429318623
$ws.Range("P13").Value = 0.06337483429318624
$ws.Range("Q13").Value = 72.64090805214389
$ws.Range("R13").Value = 653.768172469295
$ws.Range("S13").Value = 0.004477977553684103
$ws.Range("T13").Value = 0.004477977553684105
$ws.Range("G14").Value = 54.59966466666666
$ws.Range("H14").Value = 163.798994
$ws.Range("I14").Value = 0.158040016585813
$ws.Range("J14").Value = 0.158040016585813
$ws.Range("M14").Value = 12.23245266666667
$ws.Range("N14").Value = 36.697358
$ws.Range("O14").Value = 0.2605176191876535
$ws.Range("P14").Value = 0.2605176191876535
$ws.Range("Q14").Value = 667.8878136508724
$ws.Range("R14").Value = 6010.990322857852
$ws.Range("S14").Value = 0.04117220885731328
$ws.Range("T14").Value = 0.04117220885731328
$ws.Range("G15").Value = 54.59966466666666
$ws.Range("H15").Value = 163.798994
$ws.Range("I15").Value = 0.158040016585813
$ws.Range("J15").Value = 0.158040016585813
$ws.Range("O15").Value = 0.2463118768921188
$ws.Range("P15").Value = 0.2463118768921188
$ws.Range("Q15").Value = 631.4686179256948
$ws.Range("R15").Value = 5683.217561331254
$ws.Range("S15").Value = 0.03892713310931319
$ws.Range("T15").Value = 0.03892713310931319
$ws.Range("G16").Value = 54.59966466666666
$ws.Range("H16").Value = 163.798994
$ws.Range("I16").Value = 0.158040016585813
$ws.Range("J16").Value = 0.158040016585813
$ws.Range("M16").Value = 20.18080466666667
$ws.Range("N16").Value = 60.54241400000001
$ws.Range("O16").Value = 0.4297956696270414
$ws.Range("P16").Value = 0.4297956696270414
$ws.Range("Q16").Value = 1101.865167503502
$ws.Range("R16").Value = 9916.786507531517
$ws.Range("S16").Value = 0.06792491475636825
$ws.Range("T16").Value = 0.06792491475636823
$ws.Range("G17").Value = 54.59966466666666
$ws.Range("H17").Value = 163.798994
$ws.Range("I17").Value = 0.158040016585813
$ws.Range("J17").Value = 0.158040016585813
$ws.Range("M17").Value = 2.975728333333333
$ws.Range("N17").Value = 8.927185
$ws.Range("O17").Value = 0.06337483429318623
$ws.Range("P17").Value = 0.06337483429318624
$ws.Range("Q17").Value = 162.4737691390989
$ws.Range("R17").Value = 1462.26392225189
$ws.Range("S17").Value = 0.01001575986281831
$ws.Range("T17").Value = 0.01001575986281831
